$wb = $excel.ActiveWorkbook

# --- List1 sheet: update measured reading A59 (date/time value) ---
$ws1 = $wb.Worksheets.Item("List1")
$ws1.Range("A59").Value = 45538.631944444445

# --- "Měření aktivity" sheet: fill in the new measurement row 29 ---
$ws2 = $wb.Worksheets.Item("Měření aktivity")
$ws2.Range("A29").Value = 45538.633333333331
$ws2.Range("B29").Value = 0.006
$ws2.Range("C29").Value = 0.005
$ws2.Range("D29").Value = 0.005
$ws2.Range("E29").Value = 0.004
$ws2.Range("F29").Value = 0.004
$ws2.Range("G29").Value = 5.083
$ws2.Range("H29").Value = 5.076
$ws2.Range("I29").Value = 5.077
$ws2.Range("J29").Value = 5.07
$ws2.Range("K29").Value = 5.071
$ws2.Range("L29").Value = 5.072
$ws2.Range("M29").Value = 5.074
$ws2.Range("N29").Value = 5.073
$ws2.Range("O29").Value = 5.07
$ws2.Range("P29").Value = 5.074

# --- Restore/update the view state (scroll position + selection) on both sheets ---
$ws1.Activate()
$ws1.Range("A60").Select()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1

$ws2.Activate()
$ws2.Range("A29").Select()
